# OE-9291 - Make description and group required
#
# This settings-metadata spreadsheet lists, per setting "key", which
# "group" it belongs to (column D) along with its "description" (column E).
# A number of rows were mis-grouped as "Core" when they should really be
# "System" (or, for one row, "Examination"). Two description strings also
# had typos fixed. Finally, a block of rows that had been left hidden are
# now shown again.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Column D ("group") corrections: "Core" -> "System"
# ---------------------------------------------------------------------
$rowsToSystem = @(15, 44, 50, 60, 76, 83, 84, 91, 93, 101, 118, 119, 120, 121, `
    130, 131, 132, 133, 134, 135, 136, 137, 138, 139, 140, 173)

foreach ($r in $rowsToSystem) {
    $ws.Cells.Item($r, 4).Value = "System"
}

# ---------------------------------------------------------------------
# 2. Column D ("group") correction: "Core" -> "Examination"
# ---------------------------------------------------------------------
$ws.Cells.Item(172, 4).Value = "Examination"

# ---------------------------------------------------------------------
# 3. Description text fixes
# ---------------------------------------------------------------------
$ws.Cells.Item(157, 5).Value = "The body text for emails sent to Local Authorities for Certificate of Visual Impairment applications. Note that the actual application is attached to the email as a PDF"

$ws.Cells.Item(177, 5).Value = "When On, users must sign CVI events using their PIN. When Off, CVI events will be automatically signed for the current user on save"
$ws.Rows.Item(177).RowHeight = 18.75

# ---------------------------------------------------------------------
# 4. Unhide the previously hidden rows 174-181, and restore row 181's
#    height to match the other (visible) rows
# ---------------------------------------------------------------------
for ($r = 174; $r -le 181; $r++) {
    $ws.Rows.Item($r).Hidden = $false
}
$ws.Rows.Item(181).RowHeight = 19.5
